# Apply the "Add files via upload" edit: insert four new account rows into
# the "Export" sheet, keeping the existing descending sort by Saldo (column C).
# New rows (Conta, Nome, Saldo):
#   005642649, VR,         500000      -> before original row 4  (MARIANA / 004525587)
#   005338054, ELAINE,     1058.99     -> before original row 20 (JULIANA / 004813088)
#   001761119, BLUEMETRIX, 51.83       -> before original row 257 (FELIPE  / 004400640)
#   004479965, DIEGO,      48.17       -> before original row 263 (RODRIGO / 005152037)
#
# Insertions are performed from the bottom of the sheet upward so that the
# row numbers noted above (taken from the original, unmodified sheet) stay
# valid for each subsequent insertion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-SaldoRow {
    param(
        [int]$RowIndex,
        [string]$Conta,
        [string]$Nome,
        [double]$Saldo
    )

    $ws.Rows($RowIndex).Insert()

    $contaCell = $ws.Cells.Item($RowIndex, 1)
    # Force text storage so leading zeros in the account number are kept.
    $contaCell.NumberFormat = "@"
    $contaCell.Value = $Conta

    $nomeCell = $ws.Cells.Item($RowIndex, 2)
    $nomeCell.NumberFormat = "@"
    $nomeCell.Value = $Nome

    $ws.Cells.Item($RowIndex, 3).Value = $Saldo
}

# Insert from bottom to top (positional args; named args are not reliably
# bound by this runtime's PowerShell-subset parser).
Add-SaldoRow 263 "004479965" "DIEGO"      48.17
Add-SaldoRow 257 "001761119" "BLUEMETRIX" 51.83
Add-SaldoRow 20  "005338054" "ELAINE"     1058.99
Add-SaldoRow 4   "005642649" "VR"         500000
